$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = "Musical"
$ws.Range("B3").Value  = "Film-Noir"
$ws.Range("B4").Value  = "Thriller"
$ws.Range("B5").Value  = "Horror"
$ws.Range("B6").Value  = "Action"
$ws.Range("B7").Value  = "Sci-Fi"
$ws.Range("B9").Value  = "Documentary"
$ws.Range("B10").Value = "Crime"
$ws.Range("B11").Value = "Drama"
$ws.Range("B12").Value = "Western"
$ws.Range("B13").Value = "Adventure"
$ws.Range("B14").Value = "War"
$ws.Range("B16").Value = "Romance"
$ws.Range("B17").Value = "Mystery"
$ws.Range("B18").Value = "Comedy"
$ws.Range("B19").Value = "IMAX"
$ws.Range("B20").Value = "Fantasy"
